# Generate Report for Handback
# Updates the handoff/handback timestamps for the "aaf59e7b-..." row (row 2)
# across the zh-cn and de-de worksheets, and the mirrored "Latest HO Xliff
# Generate Date" on the Overview sheet (which tracks the de-de handoff time).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 is the aaf59e7b-...md source file.
# H = Latest Handoff Datetime, K = Latest Handback DateTime
$wsZhCn.Range("H2").Value = "2016-10-18 11:01:18"
$wsZhCn.Range("K2").Value = "2016-10-18 11:02:12"

# de-de sheet: row 2 is the aaf59e7b-...md source file.
$wsDeDe.Range("H2").Value = "2016-10-18 11:01:38"
$wsDeDe.Range("K2").Value = "2016-10-18 11:02:45"

# Overview sheet: row 2 mirrors the de-de handoff datetime for the same file.
$wsOverview.Range("G2").Value = "2016-10-18 11:01:38"
